$wb = $excel.ActiveWorkbook

# Labels for the new/updated rows (shared strings c1_0..c1_3)
$labels = @("c1_0", "c1_1", "c1_2", "c1_3")

# ----- Sheets with full B:K grids: grain, meat, wool -----
$gridSheets = @{
    "grain" = @(0.8517948448018292, 0.8475773720402696, 1.153146674367855, 1.147437125033402)
    "meat"  = @(0.706002939353065, 1.303792821776162, 0.6962071782238374, 1.293997060646935)
    "wool"  = @(0.706002939353065, 1.303792821776162, 0.6962071782238374, 1.293997060646935)
}

foreach ($name in $gridSheets.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $values = $gridSheets[$name]

    for ($i = 0; $i -lt 4; $i++) {
        $row = 2 + $i

        if ($row -gt 2) {
            # Copy formatting from row 2's label cell down to the new label cell
            $ws.Range("A2").Copy()
            $ws.Range("A$row").PasteSpecial(-4122)
        }

        $ws.Range("A$row").Value = $labels[$i]
        $ws.Range("B${row}:K${row}").Value = $values[$i]
    }
}

# ----- "prob" sheet: single B column -----
$ws = $wb.Worksheets.Item("prob")
$probValues = @(0.2426304458088387, 0.2573695541911614, 0.2573695541911612, 0.2426304458088387)

for ($i = 0; $i -lt 4; $i++) {
    $row = 2 + $i

    if ($row -gt 2) {
        $ws.Range("A2").Copy()
        $ws.Range("A$row").PasteSpecial(-4122)
    }

    $ws.Range("A$row").Value = $labels[$i]
    $ws.Range("B$row").Value = $probValues[$i]
}
